# Update "MarginFileClalOnlyApril2017" Rates sheet:
#  - refresh the Alt60NoTsamudIsrael block (rows 146-151) with new margin
#    values and bold-face formatting (matching the style already used by
#    the block above it, Alt60TsamudIsrael, rows 140-145)
#  - move the active selection to the area being edited (the column
#    auto-fit width change that Excel derives from the new bold text is a
#    render-time side effect, not something this exposes a setter for)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 146 (B=1) --------------------------------------------------------
$ws.Range("A146").Font.Bold = $true
$ws.Range("C146:I146").Value = 0.004
$ws.Range("J146:N146").Value = 0.0042
$ws.Range("O146:S146").Value = 0.0045
$ws.Range("T146:X146").Value = 0.0048
$ws.Range("Y146:AC146").Value = 0.0052
$ws.Range("C146:AC146").Font.Bold = $true

# --- Row 147 (B=2) --------------------------------------------------------
$ws.Range("A147").Font.Bold = $true
$ws.Range("C147:I147").Value = 0.0045
$ws.Range("J147:N147").Value = 0.0048
$ws.Range("O147:S147").Value = 0.0052
$ws.Range("T147:X147").Value = 0.0055
$ws.Range("Y147:AC147").Value = 0.0059
$ws.Range("C147:AC147").Font.Bold = $true

# --- Row 148 (B=3) --------------------------------------------------------
$ws.Range("A148").Font.Bold = $true
$ws.Range("C148:I148").Value = 0.005
$ws.Range("J148:N148").Value = 0.0054
$ws.Range("O148:S148").Value = 0.0058
$ws.Range("T148:X148").Value = 0.0062
$ws.Range("Y148:AC148").Value = 0.0068
$ws.Range("C148:AC148").Font.Bold = $true

# --- Row 149 (B=4) --------------------------------------------------------
$ws.Range("A149").Font.Bold = $true
$ws.Range("C149:I149").Value = 0.0055
$ws.Range("J149:N149").Value = 0.006
$ws.Range("O149:S149").Value = 0.0066
$ws.Range("T149:X149").Value = 0.0069
$ws.Range("Y149:AC149").Value = 0.0078
$ws.Range("C149:AC149").Font.Bold = $true

# --- Row 150 (B=5) --------------------------------------------------------
$ws.Range("A150").Font.Bold = $true
$ws.Range("C150:I150").Value = 0.006
$ws.Range("J150:N150").Value = 0.0065
$ws.Range("O150:S150").Value = 0.0071
$ws.Range("T150:X150").Value = 0.0076
$ws.Range("Y150:AC150").Value = 0.0087
$ws.Range("C150:AC150").Font.Bold = $true

# --- Row 151 (B=6) --------------------------------------------------------
$ws.Range("A151").Font.Bold = $true
$ws.Range("C151:I151").Value = 0.0065
$ws.Range("J151:N151").Value = 0.0071
$ws.Range("O151:S151").Value = 0.0077
$ws.Range("T151:X151").Value = 0.0082
$ws.Range("Y151:AC151").Value = 0.0099
$ws.Range("C151:AC151").Font.Bold = $true

# --- View / selection ------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 130
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C145").Select() | Out-Null
